$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# --- Simple global text replacements (each occurs twice, identical replacement both times) ---
Replace-All "MARIO A. VILLANUEVA" "FRANZ R. VIDA"
Replace-All "Admin Aide III" "Casual Employee"
Replace-All "TOPS detailed at Hanggang sa Kabilang Buhay Services" "Vice Mayor'S Office"
Replace-All "December 17, 1994" "June 05, 2017"
Replace-All "his optional retirement" "end of his term"
Replace-All "February 01, 2023" "July 01, 2022"
